$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 7286
    $ws.Range("F4").Value = 5499
    $ws.Range("F15").Value = 283
}
